$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff removes the entire data row for account 004329229 / GABRIEL / 9182.17
# (all following rows shift up by one). Locate that row by its account number
# so the edit is robust, then delete the whole row.
$target = $ws.Cells.Find("004329229")
if ($target -ne $null) {
    $ws.Rows.Item($target.Row).Delete()
} else {
    # Fallback: the row is known to be row 5 in the original layout.
    $ws.Rows.Item(5).Delete()
}
